$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Combine the card detail rows (2-16) into three consolidated rows (2-4),
# each formatted as a Python-style tuple string: (name, [attr1, attr2, ...])

$ws.Range("A2").Value = "('Curse of Thirst', ['{4}{B}', 'Enchantment " + [char]0x2014 + " Aura Curse', 'Enchant player', 'At the beginning of enchanted player" + [char]0x2019 + "s upkeep, Curse of Thirst deals damage to that player equal to the number of Curses attached to them.'])"

$ws.Range("A3").Value = "('Gather the Townsfolk', ['{1}{W}', 'Sorcery', 'Create two 1/1 white Human creature tokens.', 'Fateful hour " + [char]0x2014 + " If you have 5 or less life, create five of those tokens instead.'])"

$ws.Range("A4").Value = "('Nearheath Stalker', ['{4}{R}', 'Creature " + [char]0x2014 + " Vampire Rogue', 'Undying (When this creature dies, if it had no +1/+1 counters on it, return it to the battlefield under its owner" + [char]0x2019 + "s control with a +1/+1 counter on it.)', '4/1'])"

# Remove the now-unused rows 5 through 16
$ws.Rows("5:16").Delete()
